$wb = $excel.ActiveWorkbook

# ---- site_metrics ----
$ws = $wb.Worksheets.Item("site_metrics")
$ws.Range("O15").Value = 0.005779341191785079
$ws.Range("O16").Value = 0.01789475208351017
$ws.Range("O17").Value = 0.005618578486865572
$ws.Range("AK17").Value = $true
$ws.Range("AK20").Value = $true
$ws.Range("AK21").Value = $true
$ws.Range("O22").Value = 0.04698714589605456
$ws.Range("O23").Value = 0.04900963588476684
$ws.Range("AK23").Value = $true
$ws.Range("O28").Value = 0.002786193505485551
$ws.Range("AK42").Value = $true
$ws.Range("O43").Value = 0.04307847649894682
$ws.Range("O57").Value = 0.2156407608830078
$ws.Range("AK58").Value = $true
$ws.Range("AK62").Value = $true
$ws.Range("O68").Value = 0.04734500304124714
$ws.Range("O73").Value = 0.0212974678578754
$ws.Range("O78").Value = 0.00330623289296553
$ws.Range("AK78").Value = $true
$ws.Range("O80").Value = 0.002732410185299382
$ws.Range("O81").Value = 0.004649956386151363
$ws.Range("O82").Value = 0.04609438884185086
$ws.Range("O83").Value = 0.0481278083251732
$ws.Range("O86").Value = 0.0199608926356963
$ws.Range("AK88").Value = $true
$ws.Range("AK91").Value = $true
$ws.Range("O93").Value = 0.07324438939577617
$ws.Range("O94").Value = 0.007163830006606227
$ws.Range("AK96").Value = $true
$ws.Range("AK99").Value = $true
$ws.Range("O101").Value = 0.0481647123529429
$ws.Range("O103").Value = 0.04403517620970954
$ws.Range("O104").Value = 0.05533164632087837
$ws.Range("AK119").Value = $true
$ws.Range("O121").Value = 0.01162218826395471
$ws.Range("AK126").Value = $true
$ws.Range("AK128").Value = $true
$ws.Range("O131").Value = 0.1322786678075588
$ws.Range("AK132").Value = $true
$ws.Range("O133").Value = 0.02855739494976635
$ws.Range("AK133").Value = $true
$ws.Range("AK135").Value = $true
$ws.Range("AK136").Value = $true
$ws.Range("O142").Value = 0.006994705988248208

# ---- mk_duration ----
$ws = $wb.Worksheets.Item("mk_duration")
$ws.Range("M4").Value = 0.09124320211532533
$ws.Range("N4").Value = 1.688875965185925
$ws.Range("O4").Value = 0.2192118226600985
$ws.Range("P4").Value = 89
$ws.Range("Q4").Value = 2715
$ws.Range("R4").Value = 0.2739583333333333
$ws.Range("S4").Value = 2.164583333333334
$ws.Range("M26").Value = 0.9046033291427005
$ws.Range("N26").Value = 0.1198482881916595
$ws.Range("O26").Value = 0.01724137931034483
$ws.Range("P26").Value = 7
$ws.Range("Q26").Value = 2506.333333333333
$ws.Range("S26").Value = 4
$ws.Range("M29").Value = 0.05990572743327505
$ws.Range("N29").Value = 1.881486813157878
$ws.Range("O29").Value = 0.2438423645320197
$ws.Range("P29").Value = 99
$ws.Range("Q29").Value = 2713
$ws.Range("R29").Value = 0.1696969696969697
$ws.Range("S29").Value = -0.375757575757576
$ws.Range("M40").Value = 0.9334856116416619
$ws.Range("N40").Value = 0.0834602139578355
$ws.Range("O40").Value = 0.01424501424501425
$ws.Range("P40").Value = 5
$ws.Range("Q40").Value = 2297
$ws.Range("R40").Value = 0.01515151515151518
$ws.Range("S40").Value = 4.088744588744588
$ws.Range("M43").Value = 0.2886296994062758
$ws.Range("N43").Value = 1.061132478955467
$ws.Range("O43").Value = 0.1428571428571428
$ws.Range("P43").Value = 54
$ws.Range("Q43").Value = 2494.666666666667
$ws.Range("R43").Value = 0.096875
$ws.Range("S43").Value = 2.525520833333333
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.577274383745257
$ws.Range("N46").Value = -0.5573704017131537
$ws.Range("O46").Value = -0.1029411764705882
$ws.Range("P46").Value = -14
$ws.Range("Q46").Value = 544
$ws.Range("S46").Value = 8
$ws.Range("K72").Value = "no trend"
$ws.Range("L72").Value = $false
$ws.Range("M72").Value = 0.08384031020703486
$ws.Range("N72").Value = 1.728825615270013
$ws.Range("O72").Value = 0.2597402597402597
$ws.Range("P72").Value = 60
$ws.Range("Q72").Value = 1164.666666666667
$ws.Range("R72").Value = 0.1666666666666667
$ws.Range("S72").Value = 0.08333333333333348
$ws.Range("K122").Value = "no trend"
$ws.Range("L122").Value = $false
$ws.Range("M122").Value = 0.1350836263468835
$ws.Range("N122").Value = 1.49435205382276
$ws.Range("O122").Value = 0.225296442687747
$ws.Range("P122").Value = 57
$ws.Range("Q122").Value = 1404.333333333333
$ws.Range("R122").Value = 0.3015873015873016
$ws.Range("S122").Value = 3.015873015873016
$ws.Range("K127").Value = "no trend"
$ws.Range("L127").Value = $false
$ws.Range("M127").Value = 0.05112272749065738
$ws.Range("N127").Value = 1.950448151275792
$ws.Range("O127").Value = 0.3529411764705883
$ws.Range("P127").Value = 48
$ws.Range("Q127").Value = 580.6666666666666
$ws.Range("R127").Value = 0.7895833333333333
$ws.Range("S127").Value = 1.350000000000001

# ---- mk_intra_annual ----
$ws = $wb.Worksheets.Item("mk_intra_annual")
$ws.Range("M4").Value = 0.1692354188501368
$ws.Range("N4").Value = 1.374664703360094
$ws.Range("O4").Value = 0.1748768472906404
$ws.Range("P4").Value = 71
$ws.Range("Q4").Value = 2593
$ws.Range("M26").Value = 0.713284142388857
$ws.Range("N26").Value = 0.3674489795938048
$ws.Range("O26").Value = 0.04679802955665024
$ws.Range("P26").Value = 19
$ws.Range("Q26").Value = 2399.666666666667
$ws.Range("S26").Value = 1
$ws.Range("M29").Value = 0.4170480732510811
$ws.Range("N29").Value = 0.8115531676178501
$ws.Range("O29").Value = 0.1059113300492611
$ws.Range("P29").Value = 43
$ws.Range("Q29").Value = 2678.333333333333
$ws.Range("S29").Value = 2
$ws.Range("M40").Value = 1
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 2090
$ws.Range("K43").Value = "no trend"
$ws.Range("L43").Value = $false
$ws.Range("M43").Value = 0.05668784899732704
$ws.Range("N43").Value = 1.905709939635053
$ws.Range("O43").Value = 0.2513227513227513
$ws.Range("P43").Value = 95
$ws.Range("Q43").Value = 2433
$ws.Range("R43").Value = 0.08012820512820512
$ws.Range("S43").Value = -0.08173076923076916
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.6151205041013237
$ws.Range("N46").Value = -0.5027777991522042
$ws.Range("O46").Value = -0.08823529411764706
$ws.Range("P46").Value = -12
$ws.Range("Q46").Value = 478.6666666666667
$ws.Range("S46").Value = 1
$ws.Range("M72").Value = 0.7902291355909628
$ws.Range("N72").Value = 0.2660130798453453
$ws.Range("O72").Value = 0.04329004329004329
$ws.Range("P72").Value = 10
$ws.Range("Q72").Value = 1144.666666666667
$ws.Range("S72").Value = 1
$ws.Range("M122").Value = 0.0136468100021383
$ws.Range("N122").Value = 2.466428160553326
$ws.Range("O122").Value = 0.3517786561264822
$ws.Range("P122").Value = 89
$ws.Range("Q122").Value = 1273
$ws.Range("R122").Value = 0.06666666666666667
$ws.Range("S122").Value = 0.2666666666666667
$ws.Range("M127").Value = 0.9323851505432275
$ws.Range("N127").Value = -0.08484432973359157
$ws.Range("O127").Value = 0.3517786561264822
$ws.Range("P127").Value = -3
$ws.Range("Q127").Value = 555.6666666666666
$ws.Range("S127").Value = 2
